$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'301.01"
$ws.Range("E2").Formula = "'0.64%"
$ws.Range("D3").Formula = "'31.70"
$ws.Range("E3").Formula = "'0.35%"
$ws.Range("D4").Formula = "'5.105"
$ws.Range("E4").Formula = "'-1.04%"
$ws.Range("D5").Formula = "'0.07826"
$ws.Range("E5").Formula = "'-2.48%"
$ws.Range("D6").Formula = "'2.329"
$ws.Range("E6").Formula = "'-11.23%"
$ws.Range("D7").Formula = "'7.808"
$ws.Range("E7").Formula = "'-0.51%"
$ws.Range("D8").Formula = "'3.834"
$ws.Range("E8").Formula = "'0.27%"
$ws.Range("D9").Formula = "'0.9160"
$ws.Range("E9").Formula = "'0.95%"
$ws.Range("D10").Formula = "'0.1762"
$ws.Range("E10").Formula = "'1.56%"
$ws.Range("D11").Formula = "'0.07558"
$ws.Range("E11").Formula = "'4.88%"
$ws.Range("D12").Formula = "'0.09171"
$ws.Range("E12").Formula = "'14.38%"
$ws.Range("D13").Formula = "'0.03043"
$ws.Range("E13").Formula = "'0.98%"
$ws.Range("E14").Formula = "'0.53%"
$ws.Range("D15").Formula = "'0.001508"
$ws.Range("E15").Formula = "'0.90%"
$ws.Range("D16").Formula = "'0.005853"
$ws.Range("E16").Formula = "'-1.00%"
$ws.Range("D17").Formula = "'3.471"
$ws.Range("E17").Formula = "'-1.07%"
$ws.Range("E18").Formula = "'-0.38%"
$ws.Range("E19").Formula = "'-0.35%"
$ws.Range("D20").Formula = "'0.1338"
$ws.Range("E20").Formula = "'1.52%"
$ws.Range("D21").Formula = "'4.006"
$ws.Range("E21").Formula = "'-12.66%"
$ws.Range("D22").Formula = "'0.1792"
$ws.Range("E22").Formula = "'12.10%"
$ws.Range("D23").Formula = "'0.04586"
$ws.Range("E23").Formula = "'0.13%"
$ws.Range("E24").Formula = "'-0.63%"
$ws.Range("D25").Formula = "'0.004465"
$ws.Range("E25").Formula = "'0.17%"
$ws.Range("E26").Formula = "'6.03%"
$ws.Range("E27").Formula = "'-1.28%"
$ws.Range("D39").Formula = "'0.01773"
$ws.Range("E39").Formula = "'-4.63%"
$ws.Range("D40").Formula = "'0.04851"
$ws.Range("E40").Formula = "'6.75%"
$ws.Range("D41").Formula = "'0.007201"
$ws.Range("E41").Formula = "'4.26%"
$ws.Range("D42").Formula = "'0.1360"
$ws.Range("E42").Formula = "'0.91%"
$ws.Range("D43").Formula = "'0.002189"
$ws.Range("E43").Formula = "'-2.14%"
$ws.Range("D44").Formula = "'0.01027"
$ws.Range("E44").Formula = "'-1.53%"
$ws.Range("D45").Formula = "'0.00006211"
$ws.Range("E45").Formula = "'-3.56%"
$ws.Range("E46").Formula = "'0.14%"
$ws.Range("E47").Formula = "'28.97%"
$ws.Range("D48").Formula = "'1.151"
$ws.Range("E48").Formula = "'40.30%"
$ws.Range("D49").Formula = "'0.00002101"
$ws.Range("E49").Formula = "'0.14%"
$ws.Range("E50").Formula = "'0.14%"
